$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto values per diff (Tue Mar 14 18:57:48 UTC 2023 GitHub Actions refresh)
$ws.Range("D2").Value = '25.314.63'
$ws.Range("E2").Value = '  +4.34%  '

$ws.Range("D3").Value = '1.746.70'
$ws.Range("E3").Value = '  +4.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +1.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.80'
$ws.Range("E5").Value = '  +2.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.92%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3809'
$ws.Range("E7").Value = '  +2.21%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3564'
$ws.Range("E8").Value = '  +3.73%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.43'
$ws.Range("E9").Value = '  +2.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.213'
$ws.Range("E10").Value = '  +2.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07617'
$ws.Range("E11").Value = '  +4.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.007'
$ws.Range("E12").Value = '  +1.11%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.431'
$ws.Range("E13").Value = '  +5.73%  '

$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.30'
$ws.Range("E14").Value = '  +3.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.080'
$ws.Range("E15").Value = '  +4.61%  '

$ws.Range("D16").Value = '1.752.06'
$ws.Range("E16").Value = '  +4.94%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001149'
$ws.Range("E17").Value = '  +3.50%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.005'
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06746'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '85.75'
$ws.Range("E20").Value = '  +4.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.49'
$ws.Range("E21").Value = '  +5.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.480'
$ws.Range("E22").Value = '  +5.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.02'
$ws.Range("E23").Value = '  +8.19%  '

$ws.Range("D24").Value = '25.321.77'
$ws.Range("E24").Value = '  +4.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.458'
$ws.Range("E25").Value = '  +2.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.870'
$ws.Range("E26").Value = '  +7.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.82'
$ws.Range("E27").Value = '  +6.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.82'
$ws.Range("E28").Value = '  +1.33%  '

$ws.Range("D29").Value = '1.949.39'
$ws.Range("E29").Value = '  +5.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.18'
$ws.Range("E30").Value = '  +4.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.201'
$ws.Range("E31").Value = '  +21.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.087'
$ws.Range("E32").Value = '  +11.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.235'
$ws.Range("E33").Value = '  +5.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '14.25'
$ws.Range("E34").Value = '  +14.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.805'
$ws.Range("E35").Value = '  +3.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08718'
$ws.Range("E36").Value = '  +3.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.712'
$ws.Range("E37").Value = '  +6.16%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06768'
$ws.Range("E38").Value = '  +5.51%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.321'
$ws.Range("E39").Value = '  +4.17%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02469'
$ws.Range("E40").Value = '  +4.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2241'
$ws.Range("E41").Value = '  +5.60%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.283'
$ws.Range("E42").Value = '  -1.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6533'
$ws.Range("E43").Value = '  +6.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.19'
$ws.Range("E44").Value = '  +7.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.004'
$ws.Range("E45").Value = '  +0.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6241'
$ws.Range("E46").Value = '  +4.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.890'
$ws.Range("E47").Value = '  +2.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.171'
$ws.Range("E48").Value = '  +6.94%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '131.10'
$ws.Range("E49").Value = '  +3.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07390'
$ws.Range("E50").Value = '  +3.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.48'
$ws.Range("E51").Value = '  +5.52%  '
